$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# ---------------------------------------------------------------------------
# 1) Flip the FY19 diff/pct formulas (N,O columns) and the RANK.EQ order
#    argument (P column) for every data row 2:52.
# ---------------------------------------------------------------------------
$ws.Range("N2:N52").Formula = "=M2-L2"
$ws.Range("P2:P52").Formula = "=_xlfn.RANK.EQ(O2,`$O`$2:`$O`$52,1)"

# ---------------------------------------------------------------------------
# 2) New "Question 7-3" VLOOKUP block in columns F:I, rows 54-61
# ---------------------------------------------------------------------------
$ws.Range("F54").Value = "Question 7-3"
$ws.Range("F54").Style = $ws.Range("A54").Style

$ws.Range("F55").Value = "Department"
$ws.Range("G55").Value = "FY17_diff"
$ws.Range("H55").Value = "FY18_diff"
$ws.Range("I55").Value = "FY19_diff"
$ws.Range("F55:I55").Style = $ws.Range("A55").Style

$ws.Range("F56").Value = $ws.Range("A56").Value2
$ws.Range("F57").Value = $ws.Range("A57").Value2
$ws.Range("F58").Value = $ws.Range("A58").Value2
$ws.Range("F59").Value = $ws.Range("A59").Value2
$ws.Range("F60").Value = $ws.Range("A60").Value2
$ws.Range("F61").Value = $ws.Range("A61").Value2

$ws.Range("G56:G61").Formula = "=VLOOKUP(`$A56,`$A`$2:`$P`$52,MATCH(G`$55,`$1:`$1))"
$ws.Range("H56:I61").Formula = "=VLOOKUP(`$A56,`$A`$2:`$P`$52,MATCH(H`$55,`$1:`$1))"

# ---------------------------------------------------------------------------
# 3) New "Question 7-4" XLOOKUP block in columns F:I, rows 63-70
# ---------------------------------------------------------------------------
$ws.Range("F63").Value = "Question 7-4"
$ws.Range("F63").Style = $ws.Range("A63").Style

$ws.Range("F64").Value = "Department"
$ws.Range("G64").Value = "FY17_diff"
$ws.Range("H64").Value = "FY18_diff"
$ws.Range("I64").Value = "FY19_diff"
$ws.Range("F64:I64").Style = $ws.Range("A64").Style

$ws.Range("F65").Value = $ws.Range("A65").Value2
$ws.Range("F66").Value = $ws.Range("A66").Value2
$ws.Range("F67").Value = $ws.Range("A67").Value2
$ws.Range("F68").Value = $ws.Range("A68").Value2
$ws.Range("F69").Value = $ws.Range("A69").Value2
$ws.Range("F70").Value = $ws.Range("A70").Value2

$ws.Range("G65").FormulaArray = "=_xlfn.XLOOKUP(`$F65,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(G`$64,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("H65").FormulaArray = "=_xlfn.XLOOKUP(`$F65,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(H`$64,`$A`$1:`$P`$1,0)),FALSE)"
$ws.Range("I65").FormulaArray = "=_xlfn.XLOOKUP(`$F65,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(I`$64,`$A`$1:`$P`$1,0)),FALSE)"

for ($r = 66; $r -le 70; $r++) {
    $ws.Range("G$r").FormulaArray = "=_xlfn.XLOOKUP(`$F$r,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(G`$64,`$A`$1:`$P`$1,0)),FALSE)"
    $ws.Range("H$r").FormulaArray = "=_xlfn.XLOOKUP(`$F$r,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(H`$64,`$A`$1:`$P`$1,0)),FALSE)"
    $ws.Range("I$r").FormulaArray = "=_xlfn.XLOOKUP(`$F$r,`$A`$1:`$A`$52,INDEX(`$A`$1:`$P`$52,,MATCH(I`$64,`$A`$1:`$P`$1,0)),FALSE)"
}

# ---------------------------------------------------------------------------
# 4) New "Question 7-5" INDEX/MATCH block in columns F:I, rows 72-79
# ---------------------------------------------------------------------------
$ws.Range("F72").Value = "Question 7-5"
$ws.Range("F72").Style = $ws.Range("A72").Style

$ws.Range("F73").Value = "Department"
$ws.Range("G73").Value = "FY17_diff"
$ws.Range("H73").Value = "FY18_diff"
$ws.Range("I73").Value = "FY19_diff"
$ws.Range("F73:I73").Style = $ws.Range("A73").Style

$ws.Range("F74").Value = $ws.Range("A74").Value2
$ws.Range("F75").Value = $ws.Range("A75").Value2
$ws.Range("F76").Value = $ws.Range("A76").Value2
$ws.Range("F77").Value = $ws.Range("A77").Value2
$ws.Range("F78").Value = $ws.Range("A78").Value2
$ws.Range("F79").Value = $ws.Range("A79").Value2

$ws.Range("G74:G79").Formula = "=INDEX(`$A`$1:`$P`$52,MATCH(`$F74,`$A`$1:`$A`$52,0),MATCH(G`$73,`$A`$1:`$P`$1,0))"
$ws.Range("H74:I79").Formula = "=INDEX(`$A`$1:`$P`$52,MATCH(`$F74,`$A`$1:`$A`$52,0),MATCH(H`$73,`$A`$1:`$P`$1,0))"

# ---------------------------------------------------------------------------
# 5) Column F width, sheet view (scroll position/selection)
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 31.5546875

$ws.Application.ActiveWindow.ScrollRow = 52
$ws.Range("B91").Select()

# ---------------------------------------------------------------------------
# 6) Move the Question-6 bar chart down/right to clear the new F:I tables
# ---------------------------------------------------------------------------
$cht = $ws.ChartObjects(1)
$cht.Top = $ws.Range("H85").Top
$cht.Left = $ws.Range("H85").Left
$cht.Width = $ws.Range("H85:L99").Width
$cht.Height = $ws.Range("H85:L99").Height
